$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values (price/volume/hour columns)
# so Excel does not auto-convert them to Number/Percentage types.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = '282.01'
$ws.Range("E2").Value = '6.23%'
$ws.Range("G2").Value = '18'
$ws.Range("D3").Value = '26.83'
$ws.Range("E3").Value = '-0.17%'
$ws.Range("G3").Value = '18'
$ws.Range("D4").Value = '4.954'
$ws.Range("E4").Value = '5.49%'
$ws.Range("G4").Value = '18'
$ws.Range("D5").Value = '0.06420'
$ws.Range("E5").Value = '5.58%'
$ws.Range("G5").Value = '18'
$ws.Range("D6").Value = '7.001'
$ws.Range("E6").Value = '3.92%'
$ws.Range("G6").Value = '18'
$ws.Range("D7").Value = '3.351'
$ws.Range("E7").Value = '5.67%'
$ws.Range("G7").Value = '18'
$ws.Range("D8").Value = '0.8873'
$ws.Range("E8").Value = '4.32%'
$ws.Range("G8").Value = '18'
$ws.Range("D9").Value = '1.008'
$ws.Range("E9").Value = '11.28%'
$ws.Range("G9").Value = '18'
$ws.Range("E10").Value = '5.51%'
$ws.Range("G10").Value = '18'
$ws.Range("D11").Value = '0.05187'
$ws.Range("E11").Value = '3.12%'
$ws.Range("G11").Value = '18'
$ws.Range("D12").Value = '0.07411'
$ws.Range("E12").Value = '4.30%'
$ws.Range("G12").Value = '18'
$ws.Range("D13").Value = '0.03102'
$ws.Range("E13").Value = '-1.12%'
$ws.Range("G13").Value = '18'
$ws.Range("D14").Value = '0.09052'
$ws.Range("E14").Value = '0.44%'
$ws.Range("G14").Value = '18'
$ws.Range("D15").Value = '0.001562'
$ws.Range("E15").Value = '2.05%'
$ws.Range("G15").Value = '18'
$ws.Range("D16").Value = '0.0006320'
$ws.Range("E16").Value = '4.07%'
$ws.Range("G16").Value = '18'
$ws.Range("D17").Value = '0.006003'
$ws.Range("E17").Value = '0.48%'
$ws.Range("G17").Value = '18'
$ws.Range("D18").Value = '3.496'
$ws.Range("E18").Value = '1.20%'
$ws.Range("G18").Value = '18'
$ws.Range("E19").Value = '4.69%'
$ws.Range("G19").Value = '18'
$ws.Range("E20").Value = '1.41%'
$ws.Range("G20").Value = '18'
$ws.Range("E21").Value = '3.75%'
$ws.Range("G21").Value = '18'
$ws.Range("D22").Value = '3.922'
$ws.Range("E22").Value = '-4.16%'
$ws.Range("G22").Value = '18'
$ws.Range("D23").Value = '0.04348'
$ws.Range("E23").Value = '2.19%'
$ws.Range("G23").Value = '18'
$ws.Range("D24").Value = '0.001177'
$ws.Range("E24").Value = '-0.26%'
$ws.Range("G24").Value = '18'
$ws.Range("D25").Value = '0.003695'
$ws.Range("E25").Value = '-10.55%'
$ws.Range("G25").Value = '18'
$ws.Range("E26").Value = '-0.13%'
$ws.Range("G26").Value = '18'
$ws.Range("D27").Value = '0.0001693'
$ws.Range("E27").Value = '0.64%'
$ws.Range("G27").Value = '18'
$ws.Range("G28").Value = '18'
$ws.Range("G29").Value = '18'
$ws.Range("G30").Value = '18'
$ws.Range("G31").Value = '18'
$ws.Range("G32").Value = '18'
$ws.Range("G33").Value = '18'
$ws.Range("G34").Value = '18'
$ws.Range("G35").Value = '18'
$ws.Range("G36").Value = '18'
$ws.Range("G37").Value = '18'
$ws.Range("G38").Value = '18'
$ws.Range("G39").Value = '18'
$ws.Range("D40").Value = '0.04147'
$ws.Range("E40").Value = '5.85%'
$ws.Range("G40").Value = '18'
$ws.Range("D41").Value = '0.006639'
$ws.Range("E41").Value = '59.01%'
$ws.Range("G41").Value = '18'
$ws.Range("D42").Value = '0.1181'
$ws.Range("E42").Value = '6.05%'
$ws.Range("G42").Value = '18'
$ws.Range("D43").Value = '0.002359'
$ws.Range("E43").Value = '11.71%'
$ws.Range("G43").Value = '18'
$ws.Range("D44").Value = '0.01260'
$ws.Range("E44").Value = '9.75%'
$ws.Range("G44").Value = '18'
$ws.Range("D45").Value = '0.00005261'
$ws.Range("E45").Value = '2.93%'
$ws.Range("G45").Value = '18'
$ws.Range("E46").Value = '-0.08%'
$ws.Range("G46").Value = '18'
$ws.Range("E47").Value = '1,324.71%'
$ws.Range("G47").Value = '18'
$ws.Range("D48").Value = '0.02248'
$ws.Range("E48").Value = '-8.14%'
$ws.Range("G48").Value = '18'
$ws.Range("E49").Value = '-0.08%'
$ws.Range("G49").Value = '18'
$ws.Range("E50").Value = '-0.15%'
$ws.Range("G50").Value = '18'
$ws.Range("G51").Value = '18'
